$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.556.14"
$ws.Range("E2").Value = "  +6.05%  "
$ws.Range("D3").Value = "1.819.36"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "345.27"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3827"
$ws.Range("E7").Value = "  +3.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3522"
$ws.Range("E8").Value = "  +5.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.80"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  +4.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07802"
$ws.Range("E11").Value = "  +4.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.31"
$ws.Range("E13").Value = "  +10.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.625"
$ws.Range("E14").Value = "  +6.15%  "
$ws.Range("D15").Value = "1.817.73"
$ws.Range("E15").Value = "  +5.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.232"
$ws.Range("E16").Value = "  +4.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001127"
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06738"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.32"
$ws.Range("E19").Value = "  +5.24%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.74"
$ws.Range("E21").Value = "  +8.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.537"
$ws.Range("E22").Value = "  +7.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.24"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").Value = "27.544.05"
$ws.Range("E24").Value = "  +6.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.458"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.685"
$ws.Range("E26").Value = "  +8.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.19"
$ws.Range("E27").Value = "  +15.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.504"
$ws.Range("E28").Value = "  +15.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.44"
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("D30").Value = "2.020.11"
$ws.Range("E30").Value = "  +5.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.30"
$ws.Range("E31").Value = "  +5.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.389"
$ws.Range("E32").Value = "  +7.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.079"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "14.04"
$ws.Range("E34").Value = "  +8.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08777"
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.638"
$ws.Range("E37").Value = "  +5.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.7099"
$ws.Range("E38").Value = "  +15.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2283"
$ws.Range("E39").Value = "  +6.78%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02427"
$ws.Range("E40").Value = "  +6.00%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06546"
$ws.Range("E41").Value = "  +5.23%  "
$ws.Range("E42").Value = "  +5.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.290"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.84"
$ws.Range("E44").Value = "  +2.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6633"
$ws.Range("E45").Value = "  +13.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.970"
$ws.Range("E47").Value = "  +3.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.185"
$ws.Range("E48").Value = "  +8.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.80"
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07368"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.73"
$ws.Range("E51").Value = "  +4.78%  "
